# Translate the Dutch fuel-type labels in column A (rows 2-17) to English on
# both worksheets ("Nieuw" and "Tweedehands"), and restore the sheet
# selection / active-tab state recorded in the target file.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Nieuw"
$ws2 = $wb.Worksheets.Item(2)   # "Tweedehands"

# --- Sheet "Tweedehands": translate fuel-type names (A2:A17) ---
$ws2.Range("A3").Value  = "Petrol"
$ws2.Range("A4").Value  = "Diesel"
$ws2.Range("A5").Value  = "Gas + petrol"
$ws2.Range("A6").Value  = "Electric"
$ws2.Range("A2").Value  = "Unknown"
$ws2.Range("A7").Value  = "Alternative"
$ws2.Range("A8").Value  = "Natural gas"
$ws2.Range("A9").Value  = "Petrol + electric"
$ws2.Range("A10").Value = "Diesel+electric"
$ws2.Range("A11").Value = "mengsmering"
$ws2.Range("A12").Value = "Hydrogen"
$ws2.Range("A13").Value = "Diesel + gas"
$ws2.Range("A14").Value = "Electric + liquefied petroleum gas"
$ws2.Range("A15").Value = "Bio-ethanol"
$ws2.Range("A16").Value = "Hydrogen + electric"
$ws2.Range("A17").Value = "Petrol + natural gas"

# --- Sheet "Nieuw": translate fuel-type names (A2:A17) ---
$ws1.Range("A3").Value  = "Petrol"
$ws1.Range("A4").Value  = "Diesel"
$ws1.Range("A5").Value  = "Gas + petrol"
$ws1.Range("A6").Value  = "Electric"
$ws1.Range("A2").Value  = "Unknown"
$ws1.Range("A7").Value  = "Alternative"
$ws1.Range("A8").Value  = "Natural gas"
$ws1.Range("A9").Value  = "Petrol + electric"
$ws1.Range("A10").Value = "Diesel + electric"
$ws1.Range("A11").Value = "mengsmering"
$ws1.Range("A12").Value = "Hydrogen"
$ws1.Range("A13").Value = "Diesel + gas"
$ws1.Range("A14").Value = "Electric + liquefied petroleum gas"
$ws1.Range("A15").Value = "Bio-ethanol"
$ws1.Range("A16").Value = "Hydrogen + electric"
$ws1.Range("A17").Value = "Petrol + natural gas"

# --- Restore view state: "Tweedehands" was active/selected before, now
#     "Nieuw" is the active tab, with different selections on each sheet. ---
$ws2.Select()
$ws2.Range("C22").Select()

$ws1.Select()
$ws1.Range("A16").Select()
